$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 14 - this shifts the existing rows 14:19 down to 15:20
# (values + formatting move with the row; the Hyperlinks collection is left untouched,
# matching the source scraper's behaviour of not re-homing old hyperlink objects).
$ws.Rows.Item(14).Insert()

# New listing data occupying the freshly inserted row 14.
$ws.Range("A14").Value = "2025-10-22 12:50:35"
$ws.Range("B14").Value = "【急募】WEB会計アプリのテストユーザーを募集します!"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Range("F14").Value = "https://www.lancers.jp/work/detail/5418565"
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = "◇アプリ"

# Refresh the "fetched at" timestamp across every data row (2:20) to the new run time.
$ws.Range("A2").Value = "2025-10-22 12:50:35"
$ws.Range("A3").Value = "2025-10-22 12:50:35"
$ws.Range("A4").Value = "2025-10-22 12:50:35"
$ws.Range("A5").Value = "2025-10-22 12:50:35"
$ws.Range("A6").Value = "2025-10-22 12:50:35"
$ws.Range("A7").Value = "2025-10-22 12:50:35"
$ws.Range("A8").Value = "2025-10-22 12:50:35"
$ws.Range("A9").Value = "2025-10-22 12:50:35"
$ws.Range("A10").Value = "2025-10-22 12:50:35"
$ws.Range("A11").Value = "2025-10-22 12:50:35"
$ws.Range("A12").Value = "2025-10-22 12:50:35"
$ws.Range("A13").Value = "2025-10-22 12:50:35"
$ws.Range("A14").Value = "2025-10-22 12:50:35"
$ws.Range("A15").Value = "2025-10-22 12:50:35"
$ws.Range("A16").Value = "2025-10-22 12:50:35"
$ws.Range("A17").Value = "2025-10-22 12:50:35"
$ws.Range("A18").Value = "2025-10-22 12:50:35"
$ws.Range("A19").Value = "2025-10-22 12:50:35"
$ws.Range("A20").Value = "2025-10-22 12:50:35"

# The scraper only wires up a hyperlink object for rows that are brand-new to the sheet;
# row 20 (beyond the old used range) is new, so it gets a freshly added hyperlink.
$ws.Hyperlinks.Add($ws.Range("F20"), "https://www.lancers.jp/work/detail/5418533")
